# Update the four email addresses in column D (rows 2-5) of Sheet1,
# shifting the sequence from Yasser008-011s to Yasser025-028.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "Yasser025@gmail.com"
$ws.Range("D3").Value2 = "Yasser026@gmail.com"
$ws.Range("D4").Value2 = "Yasser027@gmail.com"
$ws.Range("D5").Value2 = "Yasser028@gmail.com"

# Move the active selection from D4 to D5, matching the saved view state.
$ws.Range("D5").Select() | Out-Null
